# Revert "RESTORE: Recover all 973 original multi-industry template files"
# i.e. turn the AI/ML/Data-industry specific wording back into generic
# "Product" wording across every sheet of the workbook.

$wb = $excel.ActiveWorkbook

# Touching a row's OutlineLevel (set then reset) forces Excel to persist
# that row as an explicit, otherwise-empty <row r="N"/> element even when
# it has no cell content - used below to reproduce the blank spacer rows
# that appear in the target layout.
function Touch-EmptyRow($ws, $rowNum) {
    $row = $ws.Rows.Item($rowNum)
    $row.OutlineLevel = 1
    $row.OutlineLevel = 0
}

# ---------------------------------------------------------------------------
# Sheet: Resource Overview
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resource Overview")
$ws.Range("A2").Value = "Product Implementation Project"
$ws.Range("B6").Value = "Enterprise Product Implementation"
Touch-EmptyRow $ws 13
$ws.Range("A18").Value = "Product Design/Product"
$ws.Range("G18").Value = "Development, Python, Statistics"
$ws.Range("A20").Value = "Manufacturing Engineering"
$ws.Range("A22").Value = "Production Operations/Infrastructure"
$ws.Range("G23").Value = "Manufacturing, Communication"

# ---------------------------------------------------------------------------
# Sheet: Detailed Staffing Plan
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Detailed Staffing Plan")
$ws.Range("A1").Value = "DETProductLED STAFFING PLAN"
Touch-EmptyRow $ws 2

$ws.Range("B9").Value = "Lead Product Designer"
$ws.Range("C9").Value = "Product Design/Product"
$ws.Range("K9").Value = "Development, Advanced Engineering, Python"
$ws.Range("P9").Value = "Product Lead"

$ws.Range("B10").Value = "Senior Product Designer"
$ws.Range("C10").Value = "Product Design/Product"
$ws.Range("K10").Value = "Development, Statistics, R/Python"

$ws.Range("B11").Value = "Product Designer"
$ws.Range("C11").Value = "Product Design/Product"
$ws.Range("K11").Value = "Development, Python, Visualization"

$ws.Range("B12").Value = "Development Engineer"
$ws.Range("C12").Value = "Product Design/Product"
$ws.Range("K12").Value = "DevelopmentOps, Python, Cloud"

$ws.Range("B13").Value = "Junior Product Designer"
$ws.Range("C13").Value = "Product Design/Product"

$ws.Range("B18").Value = "Senior Manufacturing Engineer"
$ws.Range("C18").Value = "Manufacturing Engineering"
$ws.Range("K18").Value = "ETL, Spark, Manufacturing Systems"

$ws.Range("B19").Value = "Manufacturing Engineer"
$ws.Range("C19").Value = "Manufacturing Engineering"
$ws.Range("K19").Value = "SQL, Python, Data Production Lines"

$ws.Range("B20").Value = "Cloud Manufacturing Engineer"
$ws.Range("C20").Value = "Manufacturing Engineering"

$ws.Range("B23").Value = "Production Operations Engineer"
$ws.Range("C23").Value = "Production Operations/Infrastructure"
$ws.Range("P23").Value = "Production Operations Lead"

$ws.Range("C24").Value = "Production Operations/Infrastructure"

$ws.Range("K25").Value = "Change Management, Manufacturing"

$ws.Range("B26").Value = "Manufacturing Specialist"
$ws.Range("K26").Value = "Manufacturing Design, Facilitation"

# ---------------------------------------------------------------------------
# Sheet: Resource Timeline
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resource Timeline")
Touch-EmptyRow $ws 2
$ws.Range("B5").Value = "Lead Product Designer"
$ws.Range("B7").Value = "Senior Manufacturing Engineer"
$ws.Range("B9").Value = "Production Operations Engineer"
Touch-EmptyRow $ws 11

# ---------------------------------------------------------------------------
# Sheet: Skills Matrix
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Skills Matrix")
Touch-EmptyRow $ws 2
$ws.Range("C3").Value = "CAD/Design Tools"
$ws.Range("D3").Value = "Product Engineering"
$ws.Range("E3").Value = "Manufacturing Engineering"
$ws.Range("F3").Value = "Manufacturing Systems"
$ws.Range("J3").Value = "Production Operations"

$ws.Range("B5").Value = "Lead Product Designer"
$ws.Range("B7").Value = "Senior Manufacturing Engineer"
$ws.Range("B9").Value = "Production Operations Engineer"
Touch-EmptyRow $ws 11

# ---------------------------------------------------------------------------
# Sheet: Cost Analysis
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cost Analysis")
Touch-EmptyRow $ws 2
$ws.Range("A6").Value = "Product Design/Product"
$ws.Range("A8").Value = "Manufacturing Engineering"
$ws.Range("A10").Value = "Production Operations/Infrastructure"
Touch-EmptyRow $ws 14
Touch-EmptyRow $ws 15

# ---------------------------------------------------------------------------
# Sheet: Resource Risk Assessment
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resource Risk Assessment")
Touch-EmptyRow $ws 2
$ws.Range("B5").Value = "Team lacks required Development expertise"
$ws.Range("F5").Value = "Manufacturing and external consultants"
Touch-EmptyRow $ws 12
Touch-EmptyRow $ws 13
